# Apply the "markov case study nearly complete 2024-06-14" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three relative-risk input values (rows 4-6) ---
$ws.Range("B4").Value = 0.8
$ws.Range("B5").Value = 0.7
$ws.Range("B6").Value = 0.9

# --- Clear the old row 8 contents (p_mildmod_ban = 1-p_severe) ---
# Rows below keep their original row numbers (no shift), so just
# clear the cells rather than deleting the whole row.
$ws.Range("A8:C8").ClearContents()

# --- Add the new row at the bottom: B17 = B13/p_injury ---
$ws.Range("B17").Formula = "=B13/p_injury"

# --- Sheet view cosmetics: scroll position + selection ---
$ws.Range("A11").Select()
$ws.Application.ActiveWindow.ScrollRow = 2

$wb.Application.Calculate()
